# Updated cryptos list on Thu Mar  2 11:50:34 UTC 2023 with GitHub Actions
#
# For each coin row: column E (hourly % change) is always refreshed, and
# column D (price) is refreshed when the price moved. Prices are stored as
# plain text in this sheet (not numbers), so we briefly force a Text number
# format before writing the value and then clear the formatting again -
# this stops Excel from "helpfully" reinterpreting a numeric-looking price
# like 299.12 as a Number while leaving the cell's style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = '23.423.15'; Volume = '  -1.31%  ' }
    @{ Row = 3; Price = '1.644.64'; Volume = '  -0.70%  ' }
    @{ Row = 4; Price = '1.001'; Volume = '  -0.16%  ' }
    @{ Row = 5; Price = $null; Volume = '  -0.08%  ' }
    @{ Row = 6; Price = '299.12'; Volume = '  -1.43%  ' }
    @{ Row = 7; Price = '0.3787'; Volume = '  -0.35%  ' }
    @{ Row = 8; Price = '49.94'; Volume = '  -2.09%  ' }
    @{ Row = 10; Price = '0.08074'; Volume = '  -1.88%  ' }
    @{ Row = 11; Price = '1.213'; Volume = '  -3.37%  ' }
    @{ Row = 12; Price = '1.002'; Volume = '  +0.05%  ' }
    @{ Row = 13; Price = '22.07'; Volume = '  -2.67%  ' }
    @{ Row = 14; Price = '6.359'; Volume = '  -2.84%  ' }
    @{ Row = 15; Price = '7.318'; Volume = '  -2.19%  ' }
    @{ Row = 16; Price = '0.00001201'; Volume = '  -3.37%  ' }
    @{ Row = 17; Price = '1.646.94'; Volume = '  -0.37%  ' }
    @{ Row = 18; Price = '96.39'; Volume = '  -1.47%  ' }
    @{ Row = 19; Price = '0.06980'; Volume = '  +0.04%  ' }
    @{ Row = 20; Price = '6.719'; Volume = '  -1.30%  ' }
    @{ Row = 21; Price = '17.33'; Volume = '  -2.46%  ' }
    @{ Row = 22; Price = $null; Volume = '  -0.11%  ' }
    @{ Row = 23; Price = '12.35'; Volume = '  -3.73%  ' }
    @{ Row = 24; Price = '23.438.23'; Volume = '  -1.24%  ' }
    @{ Row = 25; Price = $null; Volume = '  -2.55%  ' }
    @{ Row = 26; Price = '2.925'; Volume = '  -4.70%  ' }
    @{ Row = 27; Price = '20.84'; Volume = '  -2.34%  ' }
    @{ Row = 28; Price = '153.20'; Volume = '  +1.00%  ' }
    @{ Row = 29; Price = '5.205'; Volume = '  -0.25%  ' }
    @{ Row = 30; Price = '132.59'; Volume = '  -1.44%  ' }
    @{ Row = 31; Price = '1.831.11'; Volume = '  -0.50%  ' }
    @{ Row = 32; Price = '6.877'; Volume = '  -0.77%  ' }
    @{ Row = 33; Price = '2.125'; Volume = '  -2.65%  ' }
    @{ Row = 34; Price = '11.39'; Volume = '  -3.56%  ' }
    @{ Row = 35; Price = '0.9810'; Volume = '  -9.04%  ' }
    @{ Row = 36; Price = '0.02701'; Volume = '  -4.34%  ' }
    @{ Row = 37; Price = '0.08736'; Volume = '  -0.94%  ' }
    @{ Row = 38; Price = '0.2428'; Volume = '  -3.89%  ' }
    @{ Row = 39; Price = '5.905'; Volume = '  -3.70%  ' }
    @{ Row = 40; Price = '0.06805'; Volume = '  -4.16%  ' }
    @{ Row = 41; Price = $null; Volume = '  -3.43%  ' }
    @{ Row = 42; Price = '0.6872'; Volume = '  -2.90%  ' }
    @{ Row = 43; Price = '1.296'; Volume = '  -3.39%  ' }
    @{ Row = 44; Price = '15.58'; Volume = '  -2.49%  ' }
    @{ Row = 45; Price = $null; Volume = '  -0.03%  ' }
    @{ Row = 46; Price = '0.6346'; Volume = '  -3.14%  ' }
    @{ Row = 47; Price = '2.253'; Volume = '  -3.40%  ' }
    @{ Row = 48; Price = '3.908'; Volume = '  -1.38%  ' }
    @{ Row = 49; Price = $null; Volume = '  -3.05%  ' }
    @{ Row = 50; Price = '127.09'; Volume = '  -0.84%  ' }
    @{ Row = 51; Price = '1.143'; Volume = '  -4.27%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D$($u.Row)")
        # Force Text format so a numeric-looking price string (e.g. "299.12")
        # is stored as text, matching the rest of the column, then drop the
        # temporary number format so the cell style is left unchanged.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.ClearFormats()
    }
    $ws.Range("E$($u.Row)").Value = $u.Volume
}
